$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking text (e.g. "1.005", "26.930.70") which
# Excel would otherwise auto-convert to a real number, losing the original text layout.
# Force those cells to Text format before writing the new values.
# (NumberFormat is applied per contiguous block since it only affects the first
# area of a multi-area/union range.)
$ws.Range("D2:D25").NumberFormat = "@"
$ws.Range("D27:D32").NumberFormat = "@"
$ws.Range("D34:D36").NumberFormat = "@"
$ws.Range("D38:D41").NumberFormat = "@"
$ws.Range("D44:D48").NumberFormat = "@"
$ws.Range("D50:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.930.70"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.813.08"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "311.38"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "0.4281"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("D8").Value = "0.3684"
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").Value = "0.07263"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "0.8607"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "2.037.35"
$ws.Range("E11").Value = "  +12.50%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "21.19"
$ws.Range("E12").Value = "  +4.59%  "
$ws.Range("D13").Value = "6.627"
$ws.Range("E13").Value = "  +4.60%  "
$ws.Range("D14").Value = "5.378"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "0.06894"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "80.60"
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("D17").Value = "1.005"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "0.000008858"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "15.18"
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("D21").Value = "26.973.17"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "5.179"
$ws.Range("E22").Value = "  +3.41%  "
$ws.Range("D23").Value = "11.01"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "2.276.90"
$ws.Range("E24").Value = "  +12.36%  "
$ws.Range("D25").Value = "153.74"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("D27").Value = "18.25"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").Value = "5.204"
$ws.Range("E28").Value = "  +3.22%  "
$ws.Range("D29").Value = "1.884"
$ws.Range("E29").Value = "  +15.68%  "
$ws.Range("D30").Value = "114.89"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").Value = "0.08933"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "0.7411"
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("E33").Value = "  +6.39%  "
$ws.Range("D34").Value = "4.425"
$ws.Range("E34").Value = "  +2.25%  "
$ws.Range("D35").Value = "2.800"
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").Value = "1.009"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  +3.46%  "
$ws.Range("D38").Value = "0.05207"
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("D39").Value = "0.01921"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").Value = "0.5076"
$ws.Range("E40").Value = "  +3.31%  "
$ws.Range("D41").Value = "2.743"
$ws.Range("E41").Value = "  +9.40%  "
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("E43").Value = "  +5.44%  "
$ws.Range("D44").Value = "8.252"
$ws.Range("E44").Value = "  +4.02%  "
$ws.Range("B45").Value = "PaxosStandard"
$ws.Range("C45").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D45").Value = "1.005"
$ws.Range("E45").Value = "  -29.94%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "106.76"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.37"
$ws.Range("E47").Value = "  +3.39%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  +5.16%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "0.4575"
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06285"
$ws.Range("E51").Value = "  +0.98%  "
